$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1346003333333333
$ws.Range("H2").Value = 0.403801
$ws.Range("I2").Value = 0.009651054304565105
$ws.Range("J2").Value = 0.009651054304565105
$ws.Range("M2").Value = 0.110552
$ws.Range("N2").Value = 0.331656
$ws.Range("O2").Value = 0.01126249561724847
$ws.Range("P2").Value = 0.01126249561724847
$ws.Range("Q2").Value = 0.01488033605066667
$ws.Range("R2").Value = 0.133923024456
$ws.Range("S2").Value = 0.0001086949568069915
$ws.Range("T2").Value = 0.0001086949568069915
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1346003333333333
$ws.Range("H3").Value = 0.403801
$ws.Range("I3").Value = 0.009651054304565105
$ws.Range("J3").Value = 0.009651054304565105
$ws.Range("O3").Value = 0.9181055646724333
$ws.Range("P3").Value = 0.9181055646724334
$ws.Range("Q3").Value = 1.213027715756889
$ws.Range("R3").Value = 10.917249441812
$ws.Range("S3").Value = 0.008860686661977065
$ws.Range("T3").Value = 0.008860686661977065
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1346003333333333
$ws.Range("H4").Value = 0.403801
$ws.Range("I4").Value = 0.009651054304565105
$ws.Range("J4").Value = 0.009651054304565105
$ws.Range("M4").Value = 0.6933189999999999
$ws.Range("N4").Value = 2.079957
$ws.Range("O4").Value = 0.07063193971031816
$ws.Range("P4").Value = 0.07063193971031817
$ws.Range("Q4").Value = 0.09332096850633333
$ws.Range("R4").Value = 0.8398887165569999
$ws.Range("S4").Value = 0.0006816726857810491
$ws.Range("T4").Value = 0.0006816726857810492
$ws.Range("I5").Value = 0.8124788779145131
$ws.Range("J5").Value = 0.8124788779145132
$ws.Range("M5").Value = 0.110552
$ws.Range("N5").Value = 0.331656
$ws.Range("O5").Value = 0.01126249561724847
$ws.Range("P5").Value = 0.01126249561724847
$ws.Range("Q5").Value = 1.252708601144
$ws.Range("R5").Value = 11.274377410296
$ws.Range("S5").Value = 0.009150539801619161
$ws.Range("T5").Value = 0.009150539801619163
$ws.Range("I6").Value = 0.8124788779145131
$ws.Range("J6").Value = 0.8124788779145132
$ws.Range("O6").Value = 0.9181055646724333
$ws.Range("P6").Value = 0.9181055646724334
$ws.Range("S6").Value = 0.7459413789921291
$ws.Range("T6").Value = 0.7459413789921293
$ws.Range("I7").Value = 0.8124788779145131
$ws.Range("J7").Value = 0.8124788779145132
$ws.Range("M7").Value = 0.6933189999999999
$ws.Range("N7").Value = 2.079957
$ws.Range("O7").Value = 0.07063193971031816
$ws.Range("P7").Value = 0.07063193971031817
$ws.Range("Q7").Value = 7.856272836643
$ws.Range("R7").Value = 70.70645552978699
$ws.Range("S7").Value = 0.05738695912076484
$ws.Range("T7").Value = 0.05738695912076486
$ws.Range("G8").Value = 2.4807
$ws.Range("H8").Value = 7.4421
$ws.Range("I8").Value = 0.1778700677809217
$ws.Range("J8").Value = 0.1778700677809217
$ws.Range("M8").Value = 0.110552
$ws.Range("N8").Value = 0.331656
$ws.Range("O8").Value = 0.01126249561724847
$ws.Range("P8").Value = 0.01126249561724847
$ws.Range("Q8").Value = 0.2742463464
$ws.Range("R8").Value = 2.4682171176
$ws.Range("S8").Value = 0.002003260858822319
$ws.Range("T8").Value = 0.00200326085882232
$ws.Range("G9").Value = 2.4807
$ws.Range("H9").Value = 7.4421
$ws.Range("I9").Value = 0.1778700677809217
$ws.Range("J9").Value = 0.1778700677809217
$ws.Range("O9").Value = 0.9181055646724333
$ws.Range("P9").Value = 0.9181055646724334
$ws.Range("Q9").Value = 22.3562437028
$ws.Range("R9").Value = 201.2061933252
$ws.Range("S9").Value = 0.1633034990183271
$ws.Range("T9").Value = 0.1633034990183271
$ws.Range("G10").Value = 2.4807
$ws.Range("H10").Value = 7.4421
$ws.Range("I10").Value = 0.1778700677809217
$ws.Range("J10").Value = 0.1778700677809217
$ws.Range("M10").Value = 0.6933189999999999
$ws.Range("N10").Value = 2.079957
$ws.Range("O10").Value = 0.07063193971031816
$ws.Range("P10").Value = 0.07063193971031817
$ws.Range("Q10").Value = 1.7199164433
$ws.Range("R10").Value = 15.4792479897
$ws.Range("S10").Value = 0.01256330790377227
$ws.Range("T10").Value = 0.01256330790377227
